$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price record was added for "Jengibre" (ginger) at the
# "Terminal La Palmera de La Serena" market. It belongs chronologically
# right after the current row 55, so insert a fresh row at position 56,
# which pushes the previous rows 56-144 down to 57-145 (matching the
# target diff where every existing record below shifts down by one row).
$ws.Rows(56).Insert()

# Populate the newly inserted row 56 with the new record's data.
$ws.Range("A56").Value = 8
$ws.Range("B56").Value = "Terminal La Palmera de La Serena"
$ws.Range("C56").Value = "Coquimbo"
$ws.Range("D56").Value = 45100
$ws.Range("E56").Value = 4
$ws.Range("F56").Value = 100114007
$ws.Range("G56").Value = "Jengibre"
$ws.Range("H56").Value = "Sin especificar"
$ws.Range("I56").Value = "Primera"
$ws.Range("J56").Value = 240
$ws.Range("K56").Value = 17000
$ws.Range("L56").Value = 18000
$ws.Range("M56").Value = 17500
$ws.Range("N56").Value = "$/caja 13 kilos"
$ws.Range("O56").Value = "Perú"
$ws.Range("P56").Value = 1346
$ws.Range("Q56").Value = 13
$ws.Range("R56").Value = "Hortaliza"
